$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in B1: "Ontological  concept ID" (double space) -> "Ontological concept ID"
$ws.Range("B1").Value = "Ontological concept ID"

# Remove the hyperlinks that live in column H before removing the column itself
$ws.Hyperlinks.Delete()

# Delete the entire "Author's email" column (H) along with its data
$ws.Columns("H").Delete()

# Update the selected cell/range to match the target state
$ws.Range("A7").Select()
